$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New upload ("subir archivo de datos") and download ("obtener archivo de
# datos") routes were implemented, so their "terminado" (finished) column
# flips from "no" to "si".
$ws.Range("I7").Value = "si"
$ws.Range("I8").Value = "si"

# Leave the view where the author ended up after making the edit: scrolled
# down a few rows with I9 as the active selection.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I9").Select()
